# Cập nhật phân công
# Adds a new "Phân công" column (F) to Table1 on Sheet1, fills in the
# assignment values for a handful of rows, widens the new column and
# updates the sheet's active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# Grow the table by one column (A3:E48 -> A3:F48); this also extends the
# AutoFilter range and adds the 6th <tableColumn>.
$newCol = $tbl.ListColumns.Add()

# Give the new header cell the same look as the other header cells
# (border + wrapped Times New Roman text), then set its caption.
$ws.Range("B3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = "Phân công"
$excel.CutCopyMode = $false

# Fill in the assignment values that were added for this change.
$ws.Range("F4").Value = "Lê, Bùi, Tú"
$ws.Range("F5").Value = "Lê"
$ws.Range("F9").Value = "Huy"
$ws.Range("F10").Value = "Huy"
$ws.Range("F21").Value = "Nhi, Huy"
$ws.Range("F22").Value = "Nhi, Huy"

# Widen the new column so the names fit.
$ws.Columns.Item(6).ColumnWidth = 11.86

# Reflect where the author was last looking in the sheet.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("F23").Select()
